$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = 84797
$ws.Range("E23").Value = "love is in the air  beef fondue   sauces"
$ws.Range("D24").Value = 109439
$ws.Range("E24").Value = "berry  good sandwich spread"
$ws.Range("D25").Value = 42522
$ws.Range("E25").Value = "the man s  taco dip"
$ws.Range("D26").Value = 62368
$ws.Range("E26").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D126").Value = 38276
$ws.Range("E126").Value = "now and later  vegetarian empanadas"
$ws.Range("D173").Value = 62368
$ws.Range("E173").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D174").Value = 39363
$ws.Range("E174").Value = "the best  banana bread  or muffins"
$ws.Range("D201").Value = 41756
$ws.Range("E201").Value = "souper  easy sweet   sour meatballs"
$ws.Range("D224").Value = 112140
$ws.Range("E224").Value = "all in the kitchen  chili"
$ws.Range("D229").Value = 44045
$ws.Range("E229").Value = "mennonite  corn fritters"
$ws.Range("D230").Value = 87098
$ws.Range("E230").Value = "homemade  vegetable soup from a can"
$ws.Range("D231").Value = 112140
$ws.Range("E231").Value = "all in the kitchen  chili"
$ws.Range("D254").Value = 23933
$ws.Range("E254").Value = "chinese  candy"
$ws.Range("D278").Value = 112140
$ws.Range("E278").Value = "all in the kitchen  chili"
$ws.Range("D281").Value = 58224
$ws.Range("E281").Value = "immoral  sandwich filling  loose meat"
$ws.Range("D300").Value = 74805
$ws.Range("E300").Value = "never weep  whipped cream"
$ws.Range("D301").Value = 52804
$ws.Range("E301").Value = "jiffy  extra moist carrot cake"
$ws.Range("D323").Value = 75452
$ws.Range("E323").Value = "beat this  banana bread"
$ws.Range("D324").Value = 83062
$ws.Range("E324").Value = "spicy  banana bread"
$ws.Range("D325").Value = 39363
$ws.Range("E325").Value = "the best  banana bread  or muffins"
$ws.Range("D326").Value = 95926
$ws.Range("E326").Value = "say what   banana sandwich"
$ws.Range("D400").Value = 59952
$ws.Range("E400").Value = "global gourmet  taco casserole"
$ws.Range("D401").Value = 44123
$ws.Range("E401").Value = "george s at the cove  black bean soup"
$ws.Range("D448").Value = 38276
$ws.Range("E448").Value = "now and later  vegetarian empanadas"
$ws.Range("D450").Value = 67888
$ws.Range("E450").Value = "backyard style  barbecued ribs"
$ws.Range("D451").Value = 64045
$ws.Range("E451").Value = "some like it hot"
$ws.Range("D473").Value = 42570
$ws.Range("E473").Value = "pick me up  party chicken kabobs"
$ws.Range("D474").Value = 58224
$ws.Range("E474").Value = "immoral  sandwich filling  loose meat"
$ws.Range("D555").Value = 30131
$ws.Range("E555").Value = "momma s special  marinade"
$ws.Range("D556").Value = 93249
$ws.Range("E556").Value = "grilled  ranch bread"
$ws.Range("D598").Value = 83133
$ws.Range("E598").Value = "stove top  bbq  beef or pork ribs"
$ws.Range("D599").Value = 64302
$ws.Range("E599").Value = "red  macaroni salad"
$ws.Range("D698").Value = 32169
$ws.Range("E698").Value = "make that chicken dance  salsa pasta"
$ws.Range("D699").Value = 53402
$ws.Range("E699").Value = "killer  lasagna"
$ws.Range("D700").Value = 94710
$ws.Range("E700").Value = "italian  fries"
$ws.Range("D701").Value = 47366
$ws.Range("E701").Value = "forgotten  minestrone"
$ws.Range("D773").Value = 59534
$ws.Range("E773").Value = "twisted american chop suey"
$ws.Range("D774").Value = 41756
$ws.Range("E774").Value = "souper  easy sweet   sour meatballs"
$ws.Range("D775").Value = 112140
$ws.Range("E775").Value = "all in the kitchen  chili"
$ws.Range("D838").Value = 52804
$ws.Range("E838").Value = "jiffy  extra moist carrot cake"
$ws.Range("D839").Value = 26995
$ws.Range("E839").Value = "keep it going  german friendship cake"
$ws.Range("D841").Value = 27087
$ws.Range("E841").Value = "get the sensation  brownies"
$ws.Range("D850").Value = 26835
$ws.Range("E850").Value = "one bowl  perfect pound cake"
$ws.Range("D851").Value = 75452
$ws.Range("E851").Value = "beat this  banana bread"
$ws.Range("D949").Value = 62368
$ws.Range("E949").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D950").Value = 38276
$ws.Range("E950").Value = "now and later  vegetarian empanadas"
$ws.Range("D951").Value = 35653
$ws.Range("E951").Value = "make it your way  shortcakes"
